$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1207.1
$ws.Range("I2").Value = 1207.1
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1207.1
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1094.1
$ws.Range("N2").ClearContents()

$ws.Range("H21").Value = 28862.285
$ws.Range("I21").Value = 40005.668
$ws.Range("J21").Value = 20504.75
$ws.Range("K21").Value = 40005.668
$ws.Range("L21").Value = 20504.75
$ws.Range("M21").Value = -39537.668
$ws.Range("N21").Value = -21440.75

$ws.Range("H23").Value = 28862.285
$ws.Range("I23").Value = 40005.668
$ws.Range("J23").Value = 20504.75
$ws.Range("K23").Value = 40005.668
$ws.Range("L23").Value = 20504.75
$ws.Range("M23").Value = -39771.668
$ws.Range("N23").Value = -20972.75

$ws.Range("H33").Value = 536.7308
$ws.Range("I33").Value = 391.3684
$ws.Range("K33").Value = 391.3684
$ws.Range("M33").Value = -162.3684

$ws.Range("H38").Value = 945.86957
$ws.Range("I38").Value = 220.94118
$ws.Range("J38").Value = 2999.8333
$ws.Range("K38").Value = 662.82354
$ws.Range("L38").Value = 8999.499899999999
$ws.Range("M38").Value = -290.82354
$ws.Range("N38").Value = -9743.499899999999

$ws.Range("H58").Value = 22061.4
$ws.Range("J58").Value = 26840
$ws.Range("L58").Value = 80520
$ws.Range("N58").Value = -80820

$ws.Range("H82").Value = 4888
$ws.Range("I82").Value = 2813.3333
$ws.Range("K82").Value = 8439.999899999999
$ws.Range("M82").Value = -8033.999899999999

$ws.Range("H85").Value = 4888
$ws.Range("I85").Value = 2813.3333
$ws.Range("K85").Value = 8439.999899999999
$ws.Range("M85").Value = -7035.999899999999

$ws.Range("H96").Value = 726.9167
$ws.Range("I96").Value = 519.3333
$ws.Range("K96").Value = 1557.9999
$ws.Range("M96").Value = -184.9999

$ws.Range("H115").Value = 1711.8462
$ws.Range("J115").Value = 2375
$ws.Range("L115").Value = 7125
$ws.Range("N115").Value = -10259

$ws.Range("H129").Value = 3624298.2
$ws.Range("I129").Value = 41668012
$ws.Range("J129").Value = 1087.4445
$ws.Range("K129").Value = 125004036
$ws.Range("L129").Value = 3262.3335
$ws.Range("M129").Value = -124999036
$ws.Range("N129").Value = -13262.3335

$ws.Range("H135").Value = 801.63635
$ws.Range("I135").Value = 618.2
$ws.Range("J135").Value = 1515
$ws.Range("K135").Value = 5563.8
$ws.Range("L135").Value = 13635
$ws.Range("M135").Value = -3028.8
$ws.Range("N135").Value = -18705

$ws.Range("H138").Value = 3862.2456
$ws.Range("J138").Value = 4861.5386
$ws.Range("L138").Value = 14584.6158
$ws.Range("N138").Value = -24864.6158

$ws.Range("H141").Value = 372387.84
$ws.Range("I141").Value = 1316.64
$ws.Range("J141").Value = 1918517.9
$ws.Range("K141").Value = 3949.92
$ws.Range("L141").Value = 5755553.699999999
$ws.Range("M141").Value = 1230.08
$ws.Range("N141").Value = -5765913.699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2746.35
$ws.Range("I32").Value = 2312.1494
$ws.Range("J32").Value = 5652.154
$ws.Range("K32").Value = 2312.1494
$ws.Range("L32").Value = 5652.154
$ws.Range("M32").Value = -2025.1494
$ws.Range("N32").Value = -6226.154

$ws.Range("H63").Value = 3662.5
$ws.Range("I63").Value = 3243.75
$ws.Range("J63").Value = 4500
$ws.Range("K63").Value = 3243.75
$ws.Range("L63").Value = 4500
$ws.Range("M63").Value = -2557.75
$ws.Range("N63").Value = -5872

$ws.Range("H66").Value = 3662.5
$ws.Range("I66").Value = 3243.75
$ws.Range("J66").Value = 4500
$ws.Range("K66").Value = 16218.75
$ws.Range("L66").Value = 22500
$ws.Range("M66").Value = -12786.75
$ws.Range("N66").Value = -29364

$ws.Range("H80").Value = 28454.75
$ws.Range("J80").Value = 28454.75
$ws.Range("L80").Value = 28454.75
$ws.Range("N80").Value = -30450.75

$ws.Range("H83").Value = 28454.75
$ws.Range("J83").Value = 28454.75
$ws.Range("L83").Value = 85364.25
$ws.Range("N83").Value = -95348.25

$ws.Range("H101").Value = 40000
$ws.Range("J101").Value = 40000
$ws.Range("L101").Value = 40000
$ws.Range("N101").Value = -46490

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18435.688
$ws.Range("I82").Value = 2773.5557
$ws.Range("J82").Value = 38572.715
$ws.Range("K82").Value = 2773.5557
$ws.Range("L82").Value = 38572.715
$ws.Range("M82").Value = -2390.5557
$ws.Range("N82").Value = -39338.715

$ws.Range("H85").Value = 18435.688
$ws.Range("I85").Value = 2773.5557
$ws.Range("J85").Value = 38572.715
$ws.Range("K85").Value = 2773.5557
$ws.Range("L85").Value = 38572.715
$ws.Range("M85").Value = -1447.5557
$ws.Range("N85").Value = -41224.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 178003.25
$ws.Range("I31").Value = 1626.0312
$ws.Range("J31").Value = 403766.1
$ws.Range("K31").Value = 1626.0312
$ws.Range("L31").Value = 403766.1
$ws.Range("M31").Value = -1331.0312
$ws.Range("N31").Value = -404356.1

$ws.Range("H34").Value = 178003.25
$ws.Range("I34").Value = 1626.0312
$ws.Range("J34").Value = 403766.1
$ws.Range("K34").Value = 1626.0312
$ws.Range("L34").Value = 403766.1
$ws.Range("M34").Value = -1424.0312
$ws.Range("N34").Value = -404170.1

$ws.Range("H132").Value = 1854.1111
$ws.Range("I132").Value = 1354.0222
$ws.Range("J132").Value = 4354.5557
$ws.Range("K132").Value = 4062.0666
$ws.Range("L132").Value = 13063.6671
$ws.Range("M132").Value = -1532.0666
$ws.Range("N132").Value = -18123.6671

$ws.Range("H134").Value = 1306.1428
$ws.Range("I134").Value = 627.8421
$ws.Range("J134").Value = 7750
$ws.Range("K134").Value = 1883.5263
$ws.Range("L134").Value = 23250
$ws.Range("M134").Value = 651.4737
$ws.Range("N134").Value = -28320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 10200.272
$ws.Range("I87").Value = 3360.6
$ws.Range("J87").Value = 15900
$ws.Range("K87").Value = 10081.8
$ws.Range("L87").Value = 47700
$ws.Range("M87").Value = -8833.799999999999
$ws.Range("N87").Value = -50196

$ws.Range("H90").Value = 10200.272
$ws.Range("I90").Value = 3360.6
$ws.Range("J90").Value = 15900
$ws.Range("K90").Value = 30245.4
$ws.Range("L90").Value = 143100
$ws.Range("M90").Value = -24005.4
$ws.Range("N90").Value = -155580

$ws.Range("H110").Value = 2512.923
$ws.Range("I110").Value = 1031.3334
$ws.Range("J110").Value = 3782.8572
$ws.Range("K110").Value = 3094.0002
$ws.Range("L110").Value = 11348.5716
$ws.Range("M110").Value = 995.9998000000001
$ws.Range("N110").Value = -19528.5716

$ws.Range("H120").Value = 12302.857
$ws.Range("I120").Value = 11020
$ws.Range("K120").Value = 33060
$ws.Range("M120").Value = -28222

$ws.Range("H131").Value = 2224.5454
$ws.Range("I131").Value = 7516.6665
$ws.Range("J131").Value = 1388.9474
$ws.Range("K131").Value = 22549.9995
$ws.Range("L131").Value = 4166.8422
$ws.Range("M131").Value = -17509.9995
$ws.Range("N131").Value = -14246.8422

$ws.Range("H138").Value = 2619.182
$ws.Range("J138").Value = 3645
$ws.Range("L138").Value = 10935
$ws.Range("N138").Value = -21215

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4416.7144
$ws.Range("I43").Value = 1819.75
$ws.Range("J43").Value = 7879.3335
$ws.Range("K43").Value = 1819.75
$ws.Range("L43").Value = 7879.3335
$ws.Range("M43").Value = -1668.75
$ws.Range("N43").Value = -8181.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 30084
$ws.Range("J51").Value = 30084
$ws.Range("L51").Value = 30084
$ws.Range("N51").Value = -31040

$ws.Range("H136").Value = 1692.8334
$ws.Range("I136").Value = 1055.4
$ws.Range("J136").Value = 4880
$ws.Range("K136").Value = 3166.2
$ws.Range("L136").Value = 14640
$ws.Range("M136").Value = -616.2000000000003
$ws.Range("N136").Value = -19740

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1028
$ws.Range("I136").Value = 666.3611
$ws.Range("K136").Value = 1999.0833
$ws.Range("M136").Value = 550.9167000000002
